$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheet "AMSIN": row 30 was missing the normal row style used by all
#    the other data rows (A,C,D,E,F,G -> style used throughout the
#    table) and its "Run Time" value needs a tiny precision fix.
# ---------------------------------------------------------------------
$wsIn = $wb.Worksheets.Item("AMSIN")

# Re-use the formatting already applied on the rest of the table (row 2)
# so row 30 ends up visually identical to every other row.
$normalStyle = $wsIn.Cells.Item(2, 1).Style

$wsIn.Cells.Item(30, 1).Style = $normalStyle
$wsIn.Cells.Item(30, 3).Style = $normalStyle
$wsIn.Cells.Item(30, 4).Style = $normalStyle
$wsIn.Cells.Item(30, 5).Style = $normalStyle
$wsIn.Cells.Item(30, 6).Style = $normalStyle
$wsIn.Cells.Item(30, 7).Style = $normalStyle

# Correct the recorded run time (Run Time column) for the last row.
$wsIn.Cells.Item(30, 2).Value = 44951.73169278936

# ---------------------------------------------------------------------
# 2) Sheet "AMS": a new sprint run ("Aadhar" flow, pay172three) needs to
#    be appended as row 28.
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# Run Date (A28) - force text so "2023-02-03" is not auto-converted to a
# date serial number, then drop the temporary number format so the cell
# keeps the default (unstyled) look used by freshly appended rows.
$cellA = $wsAms.Cells.Item(28, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2023-02-03"
$cellA.ClearFormats()

# Run Time (B28) - numeric timestamp, re-using the same date/time display
# format already used by the rest of the column (copy format from the
# row above, then set the real value).
$wsAms.Cells.Item(27, 2).Copy() | Out-Null
$wsAms.Cells.Item(28, 2).PasteSpecial(-4122) | Out-Null
$wsAms.Cells.Item(28, 2).Value = 44960.63808569634

# Sprint Name (C28)
$cellC = $wsAms.Cells.Item(28, 3)
$cellC.NumberFormat = "@"
$cellC.Value = "pay172three"
$cellC.ClearFormats()

# Total / Pass / Fail cases and time taken
$wsAms.Cells.Item(28, 4).Value = 41
$wsAms.Cells.Item(28, 5).Value = 39
$wsAms.Cells.Item(28, 6).Value = 2
$wsAms.Cells.Item(28, 7).Value = 1.03
